$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift every timestamp in column A (rows 2-97) forward by exactly one day.
for ($r = 2; $r -le 97; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value2 = $cell.Value2 + 1
}

# Update the production values (column B) for rows 27-30 per the new model run.
$ws.Cells.Item(27, 2).Value2 = 0
$ws.Cells.Item(28, 2).Value2 = 1
$ws.Cells.Item(29, 2).Value2 = 8
$ws.Cells.Item(30, 2).Value2 = 19
